$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 46
$ws.Range("I42").Value = 46
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 138
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 92
$ws.Range("N42").ClearContents()

$ws.Range("H64").Value = 10362.25
$ws.Range("I64").Value = 4499.5
$ws.Range("J64").Value = 16225
$ws.Range("K64").Value = 4499.5
$ws.Range("L64").Value = 16225
$ws.Range("M64").Value = -4251.5
$ws.Range("N64").Value = -16721

$ws.Range("H67").Value = 10362.25
$ws.Range("I67").Value = 4499.5
$ws.Range("J67").Value = 16225
$ws.Range("K67").Value = 4499.5
$ws.Range("L67").Value = 16225
$ws.Range("M67").Value = -3641.5
$ws.Range("N67").Value = -17941

$ws.Range("H69").Value = 11842.143
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 11842.143
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 35526.429
$ws.Range("N69").Value = -37274.429

$ws.Range("H72").Value = 11842.143
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 11842.143
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 106579.287
$ws.Range("N72").Value = -115315.287

$ws.Range("H118").Value = 528.6667
$ws.Range("I118").Value = 528.6667
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1586.0001
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 70.99990000000003

$ws.Range("H125").Value = 2215.5454
$ws.Range("I125").Value = 1268.3077
$ws.Range("J125").Value = 3583.7778
$ws.Range("K125").Value = 11414.7693
$ws.Range("L125").Value = 32254.0002
$ws.Range("M125").Value = -8954.7693
$ws.Range("N125").Value = -37174.00019999999

$ws.Range("H132").Value = 2178.776
$ws.Range("I132").Value = 1877.3077
$ws.Range("J132").Value = 4791.5
$ws.Range("K132").Value = 5631.9231
$ws.Range("L132").Value = 14374.5
$ws.Range("M132").Value = -3101.9231
$ws.Range("N132").Value = -19434.5

$ws.Range("H137").Value = 516656.88
$ws.Range("I137").Value = 457247.78
$ws.Range("J137").Value = 593539.25
$ws.Range("K137").Value = 1371743.34
$ws.Range("L137").Value = 1780617.75
$ws.Range("M137").Value = -1369193.34
$ws.Range("N137").Value = -1785717.75

$ws.Range("H138").Value = 5117.659
$ws.Range("I138").Value = 2996.1875
$ws.Range("J138").Value = 6329.9287
$ws.Range("K138").Value = 8988.5625
$ws.Range("L138").Value = 18989.7861
$ws.Range("M138").Value = -3848.5625
$ws.Range("N138").Value = -29269.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3417.9167
$ws.Range("I45").Value = 2078.2942
$ws.Range("J45").Value = 6671.2856
$ws.Range("K45").Value = 2078.2942
$ws.Range("L45").Value = 6671.2856
$ws.Range("M45").Value = -1701.2942
$ws.Range("N45").Value = -7425.2856

$ws.Range("H63").Value = 4543.357
$ws.Range("I63").Value = 4123.6924
$ws.Range("J63").Value = 9999
$ws.Range("K63").Value = 4123.6924
$ws.Range("L63").Value = 9999
$ws.Range("M63").Value = -3437.6924
$ws.Range("N63").Value = -11371

$ws.Range("H66").Value = 4543.357
$ws.Range("I66").Value = 4123.6924
$ws.Range("J66").Value = 9999
$ws.Range("K66").Value = 20618.462
$ws.Range("L66").Value = 49995
$ws.Range("M66").Value = -17186.462
$ws.Range("N66").Value = -56859

$ws.Range("H132").Value = 4967.2188
$ws.Range("I132").Value = 3168.5789
$ws.Range("J132").Value = 7596
$ws.Range("K132").Value = 9505.736699999999
$ws.Range("L132").Value = 22788
$ws.Range("M132").Value = -6975.736699999999
$ws.Range("N132").Value = -27848

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37038824
$ws.Range("I20").Value = 71429920
$ws.Range("J20").Value = 2263.5386
$ws.Range("K20").Value = 71429920
$ws.Range("L20").Value = 2263.5386
$ws.Range("M20").Value = -71429673
$ws.Range("N20").Value = -2757.5386

$ws.Range("H82").Value = 57879
$ws.Range("I82").Value = 11288.25
$ws.Range("J82").Value = 120000
$ws.Range("K82").Value = 11288.25
$ws.Range("L82").Value = 120000
$ws.Range("M82").Value = -10905.25
$ws.Range("N82").Value = -120766

$ws.Range("H85").Value = 57879
$ws.Range("I85").Value = 11288.25
$ws.Range("J85").Value = 120000
$ws.Range("K85").Value = 11288.25
$ws.Range("L85").Value = 120000
$ws.Range("M85").Value = -9962.25
$ws.Range("N85").Value = -122652

$ws.Range("H99").Value = 14000
$ws.Range("I99").Value = 27500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1228053.5
$ws.Range("I31").Value = 2503116.8
$ws.Range("J31").Value = 208002.8
$ws.Range("K31").Value = 2503116.8
$ws.Range("L31").Value = 208002.8
$ws.Range("M31").Value = -2502821.8
$ws.Range("N31").Value = -208592.8

$ws.Range("H34").Value = 1228053.5
$ws.Range("I34").Value = 2503116.8
$ws.Range("J34").Value = 208002.8
$ws.Range("K34").Value = 2503116.8
$ws.Range("L34").Value = 208002.8
$ws.Range("M34").Value = -2502914.8
$ws.Range("N34").Value = -208406.8

$ws.Range("H58").Value = 208578.86
$ws.Range("I58").Value = 373306.66
$ws.Range("J58").Value = 6412.909
$ws.Range("K58").Value = 373306.66
$ws.Range("L58").Value = 6412.909
$ws.Range("M58").Value = -373103.66
$ws.Range("N58").Value = -6818.909

$ws.Range("H62").Value = 3911.1428
$ws.Range("I62").Value = 2479.8
$ws.Range("J62").Value = 7489.5
$ws.Range("K62").Value = 2479.8
$ws.Range("L62").Value = 7489.5
$ws.Range("M62").Value = -1855.8
$ws.Range("N62").Value = -8737.5

$ws.Range("H65").Value = 3911.1428
$ws.Range("I65").Value = 2479.8
$ws.Range("J65").Value = 7489.5
$ws.Range("K65").Value = 12399
$ws.Range("L65").Value = 37447.5
$ws.Range("M65").Value = -9279
$ws.Range("N65").Value = -43687.5

$ws.Range("H107").Value = 810.13794
$ws.Range("I107").Value = 772.82355
$ws.Range("J107").Value = 863
$ws.Range("K107").Value = 772.82355
$ws.Range("L107").Value = 863
$ws.Range("M107").Value = 1147.17645
$ws.Range("N107").Value = -4703

$ws.Range("H132").Value = 2680.0679
$ws.Range("I132").Value = 2119.738
$ws.Range("J132").Value = 4064.4119
$ws.Range("K132").Value = 6359.214
$ws.Range("L132").Value = 12193.2357
$ws.Range("M132").Value = -3829.214
$ws.Range("N132").Value = -17253.2357

$ws.Range("H134").Value = 671441.1
$ws.Range("I134").Value = 458503.8
$ws.Range("J134").Value = 1257018.6
$ws.Range("K134").Value = 1375511.4
$ws.Range("L134").Value = 3771055.8
$ws.Range("M134").Value = -1372976.4
$ws.Range("N134").Value = -3776125.8

$ws.Range("H136").Value = 208578.86
$ws.Range("I136").Value = 373306.66
$ws.Range("J136").Value = 6412.909
$ws.Range("K136").Value = 1119919.98
$ws.Range("L136").Value = 19238.727
$ws.Range("M136").Value = -1117369.98
$ws.Range("N136").Value = -24338.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 36.357143
$ws.Range("I6").Value = 36.357143
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 109.071429
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 3.928571000000005

$ws.Range("H34").Value = 46827.78
$ws.Range("I34").Value = 279.92307
$ws.Range("J34").Value = 107340
$ws.Range("K34").Value = 839.7692099999999
$ws.Range("L34").Value = 322020
$ws.Range("M34").Value = -755.7692099999999
$ws.Range("N34").Value = -322188

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H55").Value = 23300.125
$ws.Range("I55").Value = 6500.5
$ws.Range("J55").Value = 28900
$ws.Range("K55").Value = 19501.5
$ws.Range("L55").Value = 86700
$ws.Range("M55").Value = -19324.5
$ws.Range("N55").Value = -87054

$ws.Range("H92").Value = 890.38464
$ws.Range("I92").Value = 650.1667
$ws.Range("J92").Value = 1096.2858
$ws.Range("K92").Value = 1950.5001
$ws.Range("L92").Value = 3288.8574
$ws.Range("M92").Value = -702.5001
$ws.Range("N92").Value = -5784.857400000001

$ws.Range("H136").Value = 7317
$ws.Range("I136").Value = 5981.857
$ws.Range("J136").Value = 11990
$ws.Range("K136").Value = 17945.571
$ws.Range("L136").Value = 35970
$ws.Range("M136").Value = -12845.571
$ws.Range("N136").Value = -46170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 46205.75
$ws.Range("I32").Value = 29263
$ws.Range("J32").Value = 51853.332
$ws.Range("K32").Value = 29263
$ws.Range("L32").Value = 51853.332
$ws.Range("M32").Value = -28967
$ws.Range("N32").Value = -52445.332

$ws.Range("H101").Value = 36450.715
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 36450.715
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 36450.715
$ws.Range("N101").Value = -42940.715

$ws.Range("H132").Value = 603443.4
$ws.Range("I132").Value = 911673.1
$ws.Range("J132").Value = 81823.766
$ws.Range("K132").Value = 2735019.3
$ws.Range("L132").Value = 245471.298
$ws.Range("M132").Value = -2732489.3
$ws.Range("N132").Value = -250531.298

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 939
$ws.Range("I22").Value = 880
$ws.Range("J22").Value = 948.8333
$ws.Range("K22").Value = 880
$ws.Range("L22").Value = 948.8333
$ws.Range("M22").Value = -585
$ws.Range("N22").Value = -1538.8333

$ws.Range("H27").Value = 939
$ws.Range("I27").Value = 880
$ws.Range("J27").Value = 948.8333
$ws.Range("K27").Value = 880
$ws.Range("L27").Value = 948.8333
$ws.Range("M27").Value = -773
$ws.Range("N27").Value = -1162.8333

$ws.Range("H43").Value = 1449999.9
$ws.Range("I43").Value = 349999.5
$ws.Range("J43").Value = 2000000
$ws.Range("K43").Value = 349999.5
$ws.Range("L43").Value = 2000000
$ws.Range("M43").Value = -349806.5
$ws.Range("N43").Value = -2000386

$ws.Range("H46").Value = 3629.92
$ws.Range("I46").Value = 2843.9412
$ws.Range("J46").Value = 5300.125
$ws.Range("K46").Value = 2843.9412
$ws.Range("L46").Value = 5300.125
$ws.Range("M46").Value = -2655.9412
$ws.Range("N46").Value = -5676.125

$ws.Range("H55").Value = 671.6667
$ws.Range("I55").Value = 111.95652
$ws.Range("J55").Value = 1661.9231
$ws.Range("K55").Value = 111.95652
$ws.Range("L55").Value = 1661.9231
$ws.Range("M55").Value = 61.04348
$ws.Range("N55").Value = -2007.9231

$ws.Range("H68").Value = 50771.953
$ws.Range("I68").Value = 3750.1428
$ws.Range("J68").Value = 72715.47
$ws.Range("K68").Value = 3750.1428
$ws.Range("L68").Value = 72715.47
$ws.Range("M68").Value = -3001.1428
$ws.Range("N68").Value = -74213.47

$ws.Range("H71").Value = 50771.953
$ws.Range("I71").Value = 3750.1428
$ws.Range("J71").Value = 72715.47
$ws.Range("K71").Value = 18750.714
$ws.Range("L71").Value = 363577.35
$ws.Range("M71").Value = -15006.714
$ws.Range("N71").Value = -371065.35

$ws.Range("H93").Value = 2899.6365
$ws.Range("I93").Value = 2799.8333
$ws.Range("J93").Value = 3019.4
$ws.Range("K93").Value = 2799.8333
$ws.Range("L93").Value = 3019.4
$ws.Range("M93").Value = -1551.8333
$ws.Range("N93").Value = -5515.4

$ws.Range("H132").Value = 3839.8474
$ws.Range("I132").Value = 3149.6829
$ws.Range("J132").Value = 5411.8887
$ws.Range("K132").Value = 9449.048699999999
$ws.Range("L132").Value = 16235.6661
$ws.Range("M132").Value = -6919.048699999999
$ws.Range("N132").Value = -21295.6661

$ws.Range("H136").Value = 459764.6
$ws.Range("I136").Value = 1003862.2
$ws.Range("J136").Value = 6349.9165
$ws.Range("K136").Value = 3011586.6
$ws.Range("L136").Value = 19049.7495
$ws.Range("M136").Value = -3009036.6
$ws.Range("N136").Value = -24149.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 172915
$ws.Range("I62").Value = 339997
$ws.Range("J62").Value = 5833
$ws.Range("K62").Value = 339997
$ws.Range("L62").Value = 5833
$ws.Range("M62").Value = -339373
$ws.Range("N62").Value = -7081

$ws.Range("H65").Value = 172915
$ws.Range("I65").Value = 339997
$ws.Range("J65").Value = 5833
$ws.Range("K65").Value = 1699985
$ws.Range("L65").Value = 29165
$ws.Range("M65").Value = -1696865
$ws.Range("N65").Value = -35405

$ws.Range("H96").Value = 200977.2
$ws.Range("I96").Value = 999999
$ws.Range("J96").Value = 1221.75
$ws.Range("K96").Value = 999999
$ws.Range("L96").Value = 1221.75
$ws.Range("M96").Value = -998626
$ws.Range("N96").Value = -3967.75

$ws.Range("H100").Value = 625.8
$ws.Range("I100").Value = 625.8
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1251.6
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -710.5999999999999

$ws.Range("H122").Value = 27781366
$ws.Range("I122").Value = 50001736
$ws.Range("J122").Value = 5902.625
$ws.Range("K122").Value = 150005208
$ws.Range("L122").Value = 17707.875
$ws.Range("M122").Value = -150002758
$ws.Range("N122").Value = -22607.875
